# Slight adjustment to plate reader data.
#
# Each of the 5 sheets (WT, L157D, K160T, F193D, F193F) holds a Tecan
# plate-reader export. The "Part of Plate" / "B1-G12" note is dropped and
# the per-well table is re-laid-out to show the full row range A-H (rows
# A and H have no recorded data, they're just blank placeholder rows),
# which pushes "Start Time:", the temperature reading, and the column
# header row each up by one row.
#
# Old layout                       New layout
# ----------------------------      ----------------------------
# 28: "Part of Plate" / "B1-G12"    28: "Start Time:" / <start time>
# 29: "Start Time:" / <start time>  30: <temperature>
# 31: <temperature>                 31: "<>" header row (1..12)
# 32: "<>" header row (1..12)       32: "A"  (new, no data)
# 33: "B" + data                    33: "B" + data
# 34: "C" + data                    34: "C" + data
# 35: "D" + data                    35: "D" + data
# 36: "E" + data                    36: "E" + data
# 37: "F" + data                    37: "F" + data
# 38: "G" + data                    38: "G" + data
#                                   39: "H"  (new, no data)
# 42: "End Time:" / <end time>      42: "End Time:" / <end time>  (unchanged)

$wb = $excel.ActiveWorkbook

$grayFill = 8421504   # RGB(128,128,128) - matches the workbook's existing gray/white row-label style
$whiteFont = 16777215 # RGB(255,255,255)

$sheetNames = @("WT", "L157D", "K160T", "F193D", "F193F")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # ---- capture the values that need to move before we clear anything ----
    $startVal = $ws.Cells.Item(29, 2).Value()   # B29 -> new B28
    $tempVal  = $ws.Cells.Item(31, 2).Value()   # B31 -> new B30

    $header = @()
    for ($c = 1; $c -le 13; $c++) {
        $header += , $ws.Cells.Item(32, $c).Value()   # row 32 -> new row 31
    }

    # ---- drop the old "Part of Plate" / "B1-G12" row, old Start Time row,
    #      old temperature cell, and the old header row ----
    $ws.Cells.Item(28, 1).Clear()
    $ws.Cells.Item(28, 5).Clear()
    $ws.Rows.Item(29).Clear()
    $ws.Cells.Item(31, 2).Clear()
    $ws.Rows.Item(32).Clear()

    # ---- new row 28: Start Time ----
    $ws.Cells.Item(28, 1).Value = "Start Time:"
    $ws.Cells.Item(28, 2).Value = "'" + $startVal

    # ---- new row 30: Temperature ----
    $ws.Cells.Item(30, 2).Value = $tempVal

    # ---- new row 31: column header row ("<>" then 1..12), gray/white style ----
    for ($c = 1; $c -le 13; $c++) {
        $cell = $ws.Cells.Item(31, $c)
        $cell.Value = $header[$c - 1]
        $cell.Interior.Color = $grayFill
        $cell.Font.Color = $whiteFont
    }

    # ---- new row 32: placeholder well-row "A" (no measured data) ----
    $rowA = $ws.Cells.Item(32, 1)
    $rowA.Value = "A"
    $rowA.Interior.Color = $grayFill
    $rowA.Font.Color = $whiteFont

    # ---- new row 39: placeholder well-row "H" (no measured data) ----
    $rowH = $ws.Cells.Item(39, 1)
    $rowH.Value = "H"
    $rowH.Interior.Color = $grayFill
    $rowH.Font.Color = $whiteFont

    # ---- selection: every sheet ends up with A32:A39 selected, A32 active.
    #      Do this sheet-by-sheet (selecting moves the active sheet), and we
    #      process F193F last below so it ends up the active tab, matching
    #      the workbook's original activeTab. ----
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A32:A39").Select()
}
